# Apply weekly refresh: re-permute per-record fields (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg / unidad) across rows 2-45 of the active sheet. Mercado/Producto/Categoria
# columns (A, B, C, E-K) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44932
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 3000
$ws.Range("Q2").Value = '$/bandeja 2 kilos'
$ws.Range("R2").Value = 'Provincia de Diguillín'
$ws.Range("S2").Value = 1500
$ws.Range("T2").Value = 2
$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 4000
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("Q3").Value = '$/bandeja 2 kilos'
$ws.Range("R3").Value = 'Provincia de Diguillín'
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 2
$ws.Range("D4").Value = 44972
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2500
$ws.Range("P4").Value = 2500
$ws.Range("Q4").Value = '$/bandeja 2 kilos'
$ws.Range("R4").Value = 'Provincia de Diguillín'
$ws.Range("S4").Value = 1250
$ws.Range("T4").Value = 2
$ws.Range("D5").Value = 44951
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 2800
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 2900
$ws.Range("Q5").Value = '$/bandeja 2 kilos'
$ws.Range("R5").Value = 'Provincia de Diguillín'
$ws.Range("S5").Value = 1450
$ws.Range("T5").Value = 2
$ws.Range("D6").Value = 45006
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("Q6").Value = '$/bandeja 2 kilos'
$ws.Range("R6").Value = 'Provincia de Diguillín'
$ws.Range("S6").Value = 2000
$ws.Range("T6").Value = 2
$ws.Range("D7").Value = 44974
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("Q7").Value = '$/bandeja 2 kilos'
$ws.Range("R7").Value = 'Provincia de Diguillín'
$ws.Range("S7").Value = 1500
$ws.Range("T7").Value = 2
$ws.Range("D8").Value = 44974
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("Q8").Value = '$/bandeja 2 kilos'
$ws.Range("R8").Value = 'Provincia de Diguillín'
$ws.Range("S8").Value = 1250
$ws.Range("T8").Value = 2
$ws.Range("D9").Value = 44596
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 2500
$ws.Range("O9").Value = 2700
$ws.Range("P9").Value = 2600
$ws.Range("Q9").Value = '$/bandeja 2 kilos'
$ws.Range("R9").Value = 'Provincia de Linares'
$ws.Range("S9").Value = 1300
$ws.Range("T9").Value = 2
$ws.Range("D10").Value = 44942
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("Q10").Value = '$/bandeja 2 kilos'
$ws.Range("R10").Value = 'Provincia de Diguillín'
$ws.Range("S10").Value = 1250
$ws.Range("T10").Value = 2
$ws.Range("D11").Value = 44187
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 2800
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 2900
$ws.Range("Q11").Value = '$/bandeja 2 kilos'
$ws.Range("R11").Value = 'Provincia de Linares'
$ws.Range("S11").Value = 1450
$ws.Range("T11").Value = 2
$ws.Range("D12").Value = 44187
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 65
$ws.Range("N12").Value = 1400
$ws.Range("O12").Value = 1500
$ws.Range("P12").Value = 1446
$ws.Range("Q12").Value = '$/envase 1 kilo'
$ws.Range("R12").Value = 'Provincia de Diguillín'
$ws.Range("S12").Value = 1446
$ws.Range("T12").Value = 1
$ws.Range("D13").Value = 44979
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 3000
$ws.Range("O13").Value = 3000
$ws.Range("P13").Value = 3000
$ws.Range("Q13").Value = '$/bandeja 2 kilos'
$ws.Range("R13").Value = 'Provincia de Diguillín'
$ws.Range("S13").Value = 1500
$ws.Range("T13").Value = 2
$ws.Range("D14").Value = 44979
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 2500
$ws.Range("O14").Value = 2500
$ws.Range("P14").Value = 2500
$ws.Range("Q14").Value = '$/bandeja 2 kilos'
$ws.Range("R14").Value = 'Provincia de Diguillín'
$ws.Range("S14").Value = 1250
$ws.Range("T14").Value = 2
$ws.Range("D15").Value = 44952
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 3000
$ws.Range("O15").Value = 3000
$ws.Range("P15").Value = 3000
$ws.Range("Q15").Value = '$/bandeja 2 kilos'
$ws.Range("R15").Value = 'Provincia de Diguillín'
$ws.Range("S15").Value = 1500
$ws.Range("T15").Value = 2
$ws.Range("D16").Value = 44931
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3000
$ws.Range("P16").Value = 3000
$ws.Range("Q16").Value = '$/bandeja 2 kilos'
$ws.Range("R16").Value = 'Provincia de Diguillín'
$ws.Range("S16").Value = 1500
$ws.Range("T16").Value = 2
$ws.Range("D17").Value = 44944
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 2500
$ws.Range("O17").Value = 2500
$ws.Range("P17").Value = 2500
$ws.Range("Q17").Value = '$/bandeja 2 kilos'
$ws.Range("R17").Value = 'Provincia de Diguillín'
$ws.Range("S17").Value = 1250
$ws.Range("T17").Value = 2
$ws.Range("D18").Value = 44174
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 3700
$ws.Range("O18").Value = 3800
$ws.Range("P18").Value = 3747
$ws.Range("Q18").Value = '$/bandeja 2 kilos'
$ws.Range("R18").Value = 'Provincia de Linares'
$ws.Range("S18").Value = 1874
$ws.Range("T18").Value = 2
$ws.Range("D19").Value = 44965
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 3000
$ws.Range("O19").Value = 3000
$ws.Range("P19").Value = 3000
$ws.Range("Q19").Value = '$/bandeja 2 kilos'
$ws.Range("R19").Value = 'Provincia de Diguillín'
$ws.Range("S19").Value = 1500
$ws.Range("T19").Value = 2
$ws.Range("D20").Value = 44540
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 240
$ws.Range("N20").Value = 3500
$ws.Range("O20").Value = 3800
$ws.Range("P20").Value = 3650
$ws.Range("Q20").Value = '$/bandeja 2 kilos'
$ws.Range("R20").Value = 'Región del Maule'
$ws.Range("S20").Value = 1825
$ws.Range("T20").Value = 2
$ws.Range("D21").Value = 44988
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 30
$ws.Range("N21").Value = 3000
$ws.Range("O21").Value = 3000
$ws.Range("P21").Value = 3000
$ws.Range("Q21").Value = '$/bandeja 2 kilos'
$ws.Range("R21").Value = 'Provincia de Diguillín'
$ws.Range("S21").Value = 1500
$ws.Range("T21").Value = 2
$ws.Range("D22").Value = 44988
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 2500
$ws.Range("O22").Value = 2500
$ws.Range("P22").Value = 2500
$ws.Range("Q22").Value = '$/bandeja 2 kilos'
$ws.Range("R22").Value = 'Provincia de Diguillín'
$ws.Range("S22").Value = 1250
$ws.Range("T22").Value = 2
$ws.Range("D23").Value = 44971
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("Q23").Value = '$/bandeja 2 kilos'
$ws.Range("R23").Value = 'Provincia de Diguillín'
$ws.Range("S23").Value = 1500
$ws.Range("T23").Value = 2
$ws.Range("D24").Value = 44594
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 2500
$ws.Range("O24").Value = 2800
$ws.Range("P24").Value = 2650
$ws.Range("Q24").Value = '$/bandeja 2 kilos'
$ws.Range("R24").Value = 'Provincia de Linares'
$ws.Range("S24").Value = 1325
$ws.Range("T24").Value = 2
$ws.Range("D25").Value = 44994
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 3000
$ws.Range("O25").Value = 3200
$ws.Range("P25").Value = 3100
$ws.Range("Q25").Value = '$/bandeja 2 kilos'
$ws.Range("R25").Value = 'Provincia de Diguillín'
$ws.Range("S25").Value = 1550
$ws.Range("T25").Value = 2
$ws.Range("D26").Value = 44985
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 3000
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = 3000
$ws.Range("Q26").Value = '$/bandeja 2 kilos'
$ws.Range("R26").Value = 'Provincia de Diguillín'
$ws.Range("S26").Value = 1500
$ws.Range("T26").Value = 2
$ws.Range("D27").Value = 44985
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 2500
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2500
$ws.Range("Q27").Value = '$/bandeja 2 kilos'
$ws.Range("R27").Value = 'Provincia de Diguillín'
$ws.Range("S27").Value = 1250
$ws.Range("T27").Value = 2
$ws.Range("D28").Value = 44935
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 3000
$ws.Range("O28").Value = 3000
$ws.Range("P28").Value = 3000
$ws.Range("Q28").Value = '$/bandeja 2 kilos'
$ws.Range("R28").Value = 'Provincia de Diguillín'
$ws.Range("S28").Value = 1500
$ws.Range("T28").Value = 2
$ws.Range("D29").Value = 44953
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 3000
$ws.Range("O29").Value = 3000
$ws.Range("P29").Value = 3000
$ws.Range("Q29").Value = '$/bandeja 2 kilos'
$ws.Range("R29").Value = 'Provincia de Diguillín'
$ws.Range("S29").Value = 1500
$ws.Range("T29").Value = 2
$ws.Range("D30").Value = 45009
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 30
$ws.Range("N30").Value = 4000
$ws.Range("O30").Value = 4000
$ws.Range("P30").Value = 4000
$ws.Range("Q30").Value = '$/bandeja 2 kilos'
$ws.Range("R30").Value = 'Provincia de Diguillín'
$ws.Range("S30").Value = 2000
$ws.Range("T30").Value = 2
$ws.Range("D31").Value = 44956
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 3000
$ws.Range("O31").Value = 3000
$ws.Range("P31").Value = 3000
$ws.Range("Q31").Value = '$/bandeja 2 kilos'
$ws.Range("R31").Value = 'Provincia de Diguillín'
$ws.Range("S31").Value = 1500
$ws.Range("T31").Value = 2
$ws.Range("D32").Value = 44963
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 3000
$ws.Range("P32").Value = 3000
$ws.Range("Q32").Value = '$/bandeja 2 kilos'
$ws.Range("R32").Value = 'Provincia de Diguillín'
$ws.Range("S32").Value = 1500
$ws.Range("T32").Value = 2
$ws.Range("D33").Value = 44963
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 2500
$ws.Range("O33").Value = 2500
$ws.Range("P33").Value = 2500
$ws.Range("Q33").Value = '$/bandeja 2 kilos'
$ws.Range("R33").Value = 'Provincia de Diguillín'
$ws.Range("S33").Value = 1250
$ws.Range("T33").Value = 2
$ws.Range("D34").Value = 44966
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 30
$ws.Range("N34").Value = 2500
$ws.Range("O34").Value = 2500
$ws.Range("P34").Value = 2500
$ws.Range("Q34").Value = '$/bandeja 2 kilos'
$ws.Range("R34").Value = 'Provincia de Diguillín'
$ws.Range("S34").Value = 1250
$ws.Range("T34").Value = 2
$ws.Range("D35").Value = 44937
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 2500
$ws.Range("O35").Value = 3000
$ws.Range("P35").Value = 2750
$ws.Range("Q35").Value = '$/bandeja 2 kilos'
$ws.Range("R35").Value = 'Provincia de Diguillín'
$ws.Range("S35").Value = 1375
$ws.Range("T35").Value = 2
$ws.Range("D36").Value = 44949
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 60
$ws.Range("N36").Value = 2800
$ws.Range("O36").Value = 3000
$ws.Range("P36").Value = 2900
$ws.Range("Q36").Value = '$/bandeja 2 kilos'
$ws.Range("R36").Value = 'Provincia de Diguillín'
$ws.Range("S36").Value = 1450
$ws.Range("T36").Value = 2
$ws.Range("D37").Value = 44992
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 100
$ws.Range("N37").Value = 3000
$ws.Range("O37").Value = 3000
$ws.Range("P37").Value = 3000
$ws.Range("Q37").Value = '$/bandeja 2 kilos'
$ws.Range("R37").Value = 'Provincia de Diguillín'
$ws.Range("S37").Value = 1500
$ws.Range("T37").Value = 2
$ws.Range("D38").Value = 44967
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 50
$ws.Range("N38").Value = 3000
$ws.Range("O38").Value = 3000
$ws.Range("P38").Value = 3000
$ws.Range("Q38").Value = '$/bandeja 2 kilos'
$ws.Range("R38").Value = 'Provincia de Diguillín'
$ws.Range("S38").Value = 1500
$ws.Range("T38").Value = 2
$ws.Range("D39").Value = 44967
$ws.Range("L39").Value = 'Segunda'
$ws.Range("M39").Value = 30
$ws.Range("N39").Value = 2500
$ws.Range("O39").Value = 2500
$ws.Range("P39").Value = 2500
$ws.Range("Q39").Value = '$/bandeja 2 kilos'
$ws.Range("R39").Value = 'Provincia de Diguillín'
$ws.Range("S39").Value = 1250
$ws.Range("T39").Value = 2
$ws.Range("D40").Value = 44539
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 3800
$ws.Range("O40").Value = 4000
$ws.Range("P40").Value = 3900
$ws.Range("Q40").Value = '$/bandeja 2 kilos'
$ws.Range("R40").Value = 'Región del Maule'
$ws.Range("S40").Value = 1950
$ws.Range("T40").Value = 2
$ws.Range("D41").Value = 44181
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 65
$ws.Range("N41").Value = 3600
$ws.Range("O41").Value = 3800
$ws.Range("P41").Value = 3692
$ws.Range("Q41").Value = '$/bandeja 2 kilos'
$ws.Range("R41").Value = 'Provincia de Diguillín'
$ws.Range("S41").Value = 1846
$ws.Range("T41").Value = 2
$ws.Range("D42").Value = 44181
$ws.Range("L42").Value = 'Primera'
$ws.Range("M42").Value = 80
$ws.Range("N42").Value = 1800
$ws.Range("O42").Value = 2000
$ws.Range("P42").Value = 1875
$ws.Range("Q42").Value = '$/envase 1 kilo'
$ws.Range("R42").Value = 'Provincia de Diguillín'
$ws.Range("S42").Value = 1875
$ws.Range("T42").Value = 1
$ws.Range("D43").Value = 44970
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 50
$ws.Range("N43").Value = 3000
$ws.Range("O43").Value = 3000
$ws.Range("P43").Value = 3000
$ws.Range("Q43").Value = '$/bandeja 2 kilos'
$ws.Range("R43").Value = 'Provincia de Diguillín'
$ws.Range("S43").Value = 1500
$ws.Range("T43").Value = 2
$ws.Range("D44").Value = 44970
$ws.Range("L44").Value = 'Segunda'
$ws.Range("M44").Value = 30
$ws.Range("N44").Value = 2500
$ws.Range("O44").Value = 2500
$ws.Range("P44").Value = 2500
$ws.Range("Q44").Value = '$/bandeja 2 kilos'
$ws.Range("R44").Value = 'Provincia de Diguillín'
$ws.Range("S44").Value = 1250
$ws.Range("T44").Value = 2
$ws.Range("D45").Value = 44960
$ws.Range("L45").Value = 'Segunda'
$ws.Range("M45").Value = 60
$ws.Range("N45").Value = 2500
$ws.Range("O45").Value = 2500
$ws.Range("P45").Value = 2500
$ws.Range("Q45").Value = '$/bandeja 2 kilos'
$ws.Range("R45").Value = 'Provincia de Diguillín'
$ws.Range("S45").Value = 1250
$ws.Range("T45").Value = 2
